# Split the first "First revision" list item into two list items by
# inserting a new paragraph (same ListParagraph/numbered-list formatting)
# containing "Second revision" right after it.

$d = $word.ActiveDocument

# 1) Break "First revision" into its own paragraph by inserting a
#    paragraph mark right after it. Doing this via Find/Replace (rather
#    than Range.InsertAfter) mirrors a real edit at that caret position,
#    so Word's automatic "_GoBack" bookmark follows the edit into the
#    newly created (second) paragraph, matching the target document.
$d.Content.Find.Execute("First revision", $true, $false, $false, $false, `
    $false, $true, 1, $false, "First revision`r", 2) | Out-Null

# 2) The new, second paragraph now holds the (relocated) "_GoBack"
#    bookmark at its start. Insert "Second revision" right before the
#    paragraph mark (i.e. before the bookmark), so the final run order
#    is: text run, then bookmarkStart/bookmarkEnd - exactly like typing
#    the text and leaving the cursor (and _GoBack) at the end of it.
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$endPos = $r2.End - 1
$insertionPoint = $d.Range($endPos, $endPos)
$insertionPoint.InsertBefore("Second revision")
